$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 15, shifting existing rows 15-17 down to 16-18.
$ws.Rows.Item(15).Insert()

# The newly inserted row 15 needs to be filled with the weekly data point
# (it was copied/derived from row 14's formatting by the insert, so just
# set the values explicitly for every column).
$ws.Cells.Item(15, 1).Value = 10
$ws.Cells.Item(15, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(15, 3).Value = "La Araucanía"
$ws.Cells.Item(15, 4).NumberFormat = $ws.Cells.Item(14, 4).NumberFormat
$ws.Cells.Item(15, 4).Value = 44825
$ws.Cells.Item(15, 5).Value = 9
$ws.Cells.Item(15, 6).Value = 100112036
$ws.Cells.Item(15, 7).Value = "Caigua"
$ws.Cells.Item(15, 8).Value = "Sin especificar"
$ws.Cells.Item(15, 9).Value = "Primera"
$ws.Cells.Item(15, 10).Value = 30
$ws.Cells.Item(15, 11).Value = 20000
$ws.Cells.Item(15, 12).Value = 20000
$ws.Cells.Item(15, 13).Value = 20000
$ws.Cells.Item(15, 14).Value = "$/caja 15 kilos"
$ws.Cells.Item(15, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(15, 16).Value = 1333
$ws.Cells.Item(15, 17).Value = 15
$ws.Cells.Item(15, 18).Value = "Hortaliza"
